# ATX view added to tables
# Insert two new columns (N:O) before the existing "Досье" column and
# populate their header cells with the new "ATX" / "Наш АТХ" labels,
# mirroring the author's commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank columns at N:O - everything from column N onward shifts
# right by two columns (existing headers/data end up in P:AD instead of
# N:AB), which matches the diff's column remap.
$ws.Columns("N:O").Insert()

# New columns inherit the neighbouring header style (s="4") from the
# insert-shift automatically; just set their header text.
$ws.Range("N1").Value = "ATX"
$ws.Range("O1").Value = "Наш АТХ"

# Approximate the new columns' width to roughly match column L (the
# donor column for the insert) so the sheet still looks reasonable.
$ws.Columns("N:O").ColumnWidth = 8.5

# Match the updated selection recorded in the saved view.
$ws.Range("N1").Select()
